$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Datum: 07-01-2025", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Find result: $found  start=$($rng.Start) end=$($rng.End)"

$dateStart = $rng.Start + 7   # position right after "Datum: "; "07-01-2025" begins here

# Step 1: split "07-01-2025" into "07" | "-01-2025" (clean, ": ' run untouched)
$splitPoint = $d.Range($dateStart + 2, $dateStart + 2)
$splitPoint.InsertParagraphAfter()
$markRange = $d.Range($dateStart + 2, $dateStart + 3)
$markRange.Delete()
Write-Output "After split: [$($d.Range($rng.Start, $rng.Start + 18).Text)]"

# Step 2: temporarily bold the ": " run AND the "-01-2025" run to block merges on both sides
$colonRange = $d.Range($dateStart - 2, $dateStart)
$colonRange.Bold = 1
$afterRange = $d.Range($dateStart + 2, $dateStart + 10)
Write-Output "afterRange: [$($afterRange.Text)]"
$afterRange.Bold = 1
Write-Output "After bolding both sides: [$($d.Range($rng.Start, $rng.Start + 18).Text)]"

# Step 3: edit "07" -> "17"
$zeroRange = $d.Range($dateStart, $dateStart + 2)
$zeroRange.Text = "17"
Write-Output "After replace: [$($d.Range($rng.Start, $rng.Start + 18).Text)]"

# Step 4: restore formatting on both sides
$colonRange2 = $d.Range($dateStart - 2, $dateStart)
$colonRange2.Bold = 0
$afterRange2 = $d.Range($dateStart + 2, $dateStart + 10)
$afterRange2.Bold = 0
Write-Output "Final: [$($d.Range($rng.Start, $rng.Start + 18).Text)]"
